$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the extent of the data currently on the sheet.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

# 1. Rename header-row column captions:
#    "<Name>_old" -> "<Name>_FV2310" and "<Name>_new" -> "<Name>_FV2404"
#    (the input/output format-version suffixes used by the exporter).
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2310"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2404"
        }
    }
}

# 2. Turn the data range into an Excel Table ("ListObject") so the new
#    headers double as filter/sort controls.
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
